$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 4411.1113
$ws.Range("J40").Value = 6170
$ws.Range("K40").Value = 4411.1113
$ws.Range("L40").Value = 6170
$ws.Range("M40").Value = -4236.1113
$ws.Range("N40").Value = -6520

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8816.25
$ws.Range("I62").Value = 6647.143
$ws.Range("K62").Value = 6647.143
$ws.Range("M62").Value = -6023.143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10039.6
$ws.Range("J64").Value = 11743.454
$ws.Range("L64").Value = 11743.454
$ws.Range("N64").Value = -12239.454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 8816.25
$ws.Range("I65").Value = 6647.143
$ws.Range("K65").Value = 33235.715
$ws.Range("M65").Value = -30115.715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 10039.6
$ws.Range("J67").Value = 11743.454
$ws.Range("L67").Value = 11743.454
$ws.Range("N67").Value = -13459.454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 10799.9
$ws.Range("I74").Value = 11000
$ws.Range("J74").Value = 10714.143
$ws.Range("K74").Value = 11000
$ws.Range("L74").Value = 10714.143
$ws.Range("M74").Value = -10064
$ws.Range("N74").Value = -12586.143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 10799.9
$ws.Range("I77").Value = 11000
$ws.Range("J77").Value = 10714.143
$ws.Range("K77").Value = 55000
$ws.Range("L77").Value = 53570.715
$ws.Range("M77").Value = -50320
$ws.Range("N77").Value = -62930.715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 407.84616
$ws.Range("I92").Value = 334.72726
$ws.Range("K92").Value = 334.72726
$ws.Range("M92").Value = 913.27274

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 13889.444
$ws.Range("I63").Value = 12429.286
$ws.Range("J63").Value = 19000
$ws.Range("K63").Value = 12429.286
$ws.Range("L63").Value = 19000
$ws.Range("M63").Value = -11743.286
$ws.Range("N63").Value = -20372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 13889.444
$ws.Range("I66").Value = 12429.286
$ws.Range("J66").Value = 19000
$ws.Range("K66").Value = 62146.43
$ws.Range("L66").Value = 95000
$ws.Range("M66").Value = -58714.43
$ws.Range("N66").Value = -101864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1630.8
$ws.Range("I88").Value = 1957.6666
$ws.Range("J88").Value = 1549.0834
$ws.Range("K88").Value = 1957.6666
$ws.Range("L88").Value = 1549.0834
$ws.Range("M88").Value = -1551.6666
$ws.Range("N88").Value = -2361.0834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1630.8
$ws.Range("I91").Value = 1957.6666
$ws.Range("J91").Value = 1549.0834
$ws.Range("K91").Value = 1957.6666
$ws.Range("L91").Value = 1549.0834
$ws.Range("M91").Value = -553.6666
$ws.Range("N91").Value = -4357.0834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 87290.42999999999
$ws.Range("J92").Value = 87290.42999999999
$ws.Range("L92").Value = 87290.42999999999
$ws.Range("N92").Value = -92282.42999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 55697.5
$ws.Range("J95").Value = 60837
$ws.Range("L95").Value = 60837
$ws.Range("N95").Value = -66329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4247
$ws.Range("I110").Value = 3721.1
$ws.Range("J110").Value = 6000
$ws.Range("K110").Value = 3721.1
$ws.Range("L110").Value = 6000
$ws.Range("M110").Value = -1676.1
$ws.Range("N110").Value = -10090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1811.9032
$ws.Range("I132").Value = 1316.3334
$ws.Range("J132").Value = 2852.6
$ws.Range("K132").Value = 3949.0002
$ws.Range("L132").Value = 8557.799999999999
$ws.Range("M132").Value = -1419.0002
$ws.Range("N132").Value = -13617.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 44418.617
$ws.Range("I99").Value = 28131.2
$ws.Range("K99").Value = 28131.2
$ws.Range("M99").Value = -26633.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 15065.125
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4700
$ws.Range("I105").Value = 3133.3333
$ws.Range("K105").Value = 3133.3333
$ws.Range("M105").Value = -1386.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2709.6191
$ws.Range("I58").Value = 1837.1
$ws.Range("K58").Value = 1837.1
$ws.Range("M58").Value = -1634.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2709.6191
$ws.Range("I136").Value = 1837.1
$ws.Range("K136").Value = 5511.299999999999
$ws.Range("M136").Value = -2961.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 6967.375
$ws.Range("I138").Value = 4678.25
$ws.Range("K138").Value = 14034.75
$ws.Range("M138").Value = -8894.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15332.333
$ws.Range("I113").Value = 5999.5
$ws.Range("J113").Value = 17998.857
$ws.Range("K113").Value = 5999.5
$ws.Range("L113").Value = 17998.857
$ws.Range("N113").Value = -22338.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4758.4546
$ws.Range("I122").Value = 3952.1177
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 11856.3531
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -9406.3531
$ws.Range("N122").Value = -27400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 44667.332
$ws.Range("I40").Value = 14501
$ws.Range("K40").Value = 14501
$ws.Range("M40").Value = -14365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2482.2856
$ws.Range("I46").Value = 1159.8
$ws.Range("J46").Value = 3217
$ws.Range("K46").Value = 1159.8
$ws.Range("L46").Value = 3217
$ws.Range("M46").Value = -971.8
$ws.Range("N46").Value = -3593

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 371035.34
$ws.Range("I100").Value = 371035.34
$ws.Range("K100").Value = 371035.34
$ws.Range("M100").Value = -370494.34

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7938633.5
$ws.Range("J81").Value = 22225222
$ws.Range("L81").Value = 44450444
$ws.Range("N81").Value = -44452566

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 7938633.5
$ws.Range("J84").Value = 22225222
$ws.Range("L84").Value = 222252220
$ws.Range("N84").Value = -222262828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 87662.25
$ws.Range("J86").Value = 87662.25
$ws.Range("L86").Value = 87662.25
$ws.Range("N86").Value = -89908.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value = 87662.25
$ws.Range("J89").Value = 87662.25
$ws.Range("L89").Value = 438311.25
$ws.Range("N89").Value = -449543.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 784
$ws.Range("I100").Value = 951.3333
$ws.Range("J100").Value = 616.6667
$ws.Range("K100").Value = 1902.6666
$ws.Range("L100").Value = 1233.3334
$ws.Range("M100").Value = -1361.6666
$ws.Range("N100").Value = -2315.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M113").Value = -3829.5
